$d = $word.ActiveDocument

# 1. Remove the paragraph "The following diagram displays each deliverable of the NSW Traffic Penalty Tool. "
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "The following diagram displays each deliverable of the NSW Traffic Penalty Tool\.") {
        $p.Range.Delete()
        break
    }
}

# 2. Replace the caption suffix text (keep the leading space character
#    untouched so the replacement doesn't start exactly at the run
#    boundary and inherit formatting from the preceding field run).
$d.Content.Find.Execute("- NSW Traffic Penalty Tool WBS", $true, $false, $false, $false, $false,
                         $true, 1, $false, "– NTPT Work Breakdown Structure", 2)

# 3. Insert new paragraph text before the page break run that follows the
#    "Figure 1 ... Work Breakdown Structure" caption paragraph.
$newText = "Figure 1 displays each major deliverable and their associated components for the NTPT. " + `
           "Each work breakdown structure has four levels: main deliverable, key phases, work packages and activities. A top-down approach " + `
           "has been incorporated in forming Figure 1, with the largest item and main deliverable being the data analysis and visualisation software. The key phases include initialisation, planning, testing, backend, frontend, updates and closing. " + `
           "Each phase must be completed for the project to be completed. "

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.Length -eq 2 -and [int][char]$t[0] -eq 12) {
        if ($i -gt 1) {
            $prevText = $d.Paragraphs.Item($i - 1).Range.Text
            if ($prevText -match "Work Breakdown Structure") {
                $insertRange = $p.Range.Duplicate
                $insertRange.Collapse(1)
                $insertRange.InsertBefore($newText)
                break
            }
        }
    }
}
